$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header row 1 (Spanish -> English snake_case column names)
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Capitalize lowercase Spanish articles (de/del/el/los/la -> De/Del/El/Los/La)
# in state/municipality name cells
$ws.Range("B12").Value = "Amatenango De La Frontera"
$ws.Range("B29").Value = "Marqués De Comillas"
$ws.Range("A49").Value = "Ciudad De México"
$ws.Range("B68").Value = "San Pedro Del Gallo"
$ws.Range("A70").Value = "Estado De México"
$ws.Range("B70").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B71").Value = "Almoloya De Alquisiras"
$ws.Range("B72").Value = "Almoloya De Juárez"
$ws.Range("B79").Value = "Chapa De Mota"
$ws.Range("B83").Value = "Ecatepec De Morelos"
$ws.Range("B87").Value = "Ixtapan De La Sal"
$ws.Range("B94").Value = "Naucalpan De Juárez"
$ws.Range("B97").Value = "San Felipe Del Progreso"
$ws.Range("B105").Value = "Tenango Del Valle"
$ws.Range("B108").Value = "Tlalnepantla De Baz"
$ws.Range("B111").Value = "Valle De Chalco Solidaridad"
$ws.Range("B112").Value = "Villa De Allende"
$ws.Range("B119").Value = "Apaseo El Grande"
$ws.Range("B131").Value = "Acapulco De Juárez"
$ws.Range("B134").Value = "Atoyac De Álvarez"
$ws.Range("B135").Value = "Chilpancingo De Los Bravo"
$ws.Range("B136").Value = "Coahuayutla De José María Izazaga"
$ws.Range("B140").Value = "Coyuca De Benítez"
$ws.Range("B141").Value = "Coyuca De Catalán"
$ws.Range("B142").Value = "Cuetzala Del Progreso"
$ws.Range("B143").Value = "Cutzamala De Pinzón"
$ws.Range("B148").Value = "Huitzuco De Los Figueroa"
$ws.Range("B149").Value = "Iguala De La Independencia"
$ws.Range("B150").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B157").Value = "Técpan De Galeana"
$ws.Range("B159").Value = "Tepecoacuilco De Trujano"
$ws.Range("B161").Value = "Tixtla De Guerrero"
$ws.Range("B163").Value = "Tlapa De Comonfort"
$ws.Range("B166").Value = "Atotonilco El Grande"
$ws.Range("B167").Value = "Cuautepec De Hinojosa"
$ws.Range("B170").Value = "Jacala De Ledezma"
$ws.Range("B172").Value = "Pachuca De Soto"
$ws.Range("B175").Value = "Tepehuacán De Guerrero"
$ws.Range("B176").Value = "Tula De Allende"
$ws.Range("B179").Value = "Atemajac De Brizuela"
$ws.Range("B181").Value = "Encarnación De Díaz"
$ws.Range("B183").Value = "Lagos De Moreno"
$ws.Range("B185").Value = "San Martín De Bolaños"
$ws.Range("B186").Value = "San Miguel El Alto"
$ws.Range("B187").Value = "Tamazula De Gordiano"
$ws.Range("B190").Value = "Zapotlán El Grande"
$ws.Range("B227").Value = "Puente De Ixtla"
$ws.Range("B231").Value = "Tetela Del Volcán"
$ws.Range("B232").Value = "Tlaltizapán De Zapata"
$ws.Range("B239").Value = "Amatlán De Cañas"
$ws.Range("B245").Value = "Acatlán De Pérez Figueroa"
$ws.Range("B248").Value = "Huautla De Jiménez"
$ws.Range("B249").Value = "Mariscala De Juárez"
$ws.Range("B251").Value = "Nejapa De Madero"
$ws.Range("B252").Value = "Oaxaca De Juárez"
$ws.Range("B255").Value = "San Felipe Jalapa De Díaz"
$ws.Range("B271").Value = "Zimatlán De Álvarez"
$ws.Range("B287").Value = "Huehuetlán El Chico"
$ws.Range("B291").Value = "Izúcar De Matamoros"
$ws.Range("B294").Value = "Palmar De Bravo"
$ws.Range("B302").Value = "Tecali De Herrera"
$ws.Range("B307").Value = "Tepanco De López"
$ws.Range("B308").Value = "Tepatlaxco De Hidalgo"
$ws.Range("B311").Value = "Tetela De Ocampo"
$ws.Range("B313").Value = "Tlacotepec De Benito Juárez"
$ws.Range("B329").Value = "Cadereyta De Montes"
$ws.Range("B333").Value = "Ciudad Del Maíz"
$ws.Range("B338").Value = "San Ciro De Acosta"
$ws.Range("B342").Value = "Villa De Arista"
$ws.Range("B363").Value = "Amatlán De Los Reyes"
$ws.Range("B368").Value = "Castillo De Teayo"
$ws.Range("B373").Value = "Cosamaloapan De Carpio"
$ws.Range("B375").Value = "Ignacio De La Llave"
$ws.Range("B378").Value = "Lerdo De Tejada"
$ws.Range("B379").Value = "Martínez De La Torre"
$ws.Range("B385").Value = "Paso De Ovejas"
$ws.Range("B388").Value = "Sayula De Alemán"
$ws.Range("B406").Value = "Villa De Cos"

# Remove trailing footnote/metadata rows (411-415) that are no longer needed
$ws.Rows("411:415").Delete()
